$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values (5) into previously-empty cells
$ws.Range("C6:F6").Value = 5
$ws.Range("C9:F9").Value = 5
$ws.Range("D31:E31").Value = 5

# Update the selected/active cell to G9 (matching the saved view state)
$ws.Range("G9").Select()
